# The commit adds one new weekly price record for "Tomate" (Larga vida,
# Primera) at the top of the existing block of rows (which starts at row
# 271), pushing all the subsequent rows down by one. The sheet's
# dimension grows from A1:R338 to A1:R339.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 271; this shifts rows 271:338 down
# to 272:339 and also carries the existing row formatting (e.g. the date
# number format on column D) down with them.
$ws.Rows("271").Insert()

# Populate the newly inserted row 271 with the new record.
$ws.Range("A271").Value = 7
$ws.Range("B271").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C271").Value = "Ñuble"
$ws.Range("D271").Value = 44543
$ws.Range("E271").Value = 16
$ws.Range("F271").Value = 100112020
$ws.Range("G271").Value = "Tomate"
$ws.Range("H271").Value = "Larga vida"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 800
$ws.Range("K271").Value = 6000
$ws.Range("L271").Value = 7000
$ws.Range("M271").Value = 6500
$ws.Range("N271").Value = "`$/caja 15 kilos"
$ws.Range("O271").Value = "Región del Maule"
$ws.Range("P271").Value = 433
$ws.Range("Q271").Value = 15
$ws.Range("R271").Value = "Hortaliza"
